$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to keep a literal text value (matching the source
    # workbook, where these cells are stored as inline/shared strings)
    # even when the text looks like a number (e.g. "218.32" or "0.5220").
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.157.64'
Set-TextValue 'E2' '  -2.08%  '
Set-TextValue 'D3' '1.658.45'
Set-TextValue 'E3' '  -1.88%  '
Set-TextValue 'E4' '  +0.59%  '
Set-TextValue 'D5' '218.32'
Set-TextValue 'E5' '  +0.00%  '
Set-TextValue 'D6' '0.5220'
Set-TextValue 'E6' '  -2.79%  '
Set-TextValue 'E7' '  +0.57%  '
Set-TextValue 'D8' '0.2674'
Set-TextValue 'E8' '  -0.80%  '
Set-TextValue 'D9' '0.06327'
Set-TextValue 'E9' '  -2.02%  '
Set-TextValue 'D10' '21.08'
Set-TextValue 'E10' '  -2.62%  '
Set-TextValue 'E11' '  -0.88%  '
Set-TextValue 'D12' '4.437'
Set-TextValue 'E12' '  -2.01%  '
Set-TextValue 'D13' '1.648.93'
Set-TextValue 'E13' '  -2.46%  '
Set-TextValue 'D14' '1.886.25'
Set-TextValue 'E14' '  -1.87%  '
Set-TextValue 'D15' '0.5478'
Set-TextValue 'E15' '  -3.39%  '
Set-TextValue 'D16' '0.0₅8234'
Set-TextValue 'E16' '  -3.35%  '
Set-TextValue 'D17' '64.96'
Set-TextValue 'E17' '  -2.38%  '
Set-TextValue 'D18' '26.224.55'
Set-TextValue 'E18' '  -1.92%  '
Set-TextValue 'E19' '  +0.26%  '
Set-TextValue 'D20' '4.668'
Set-TextValue 'E20' '  -3.40%  '
Set-TextValue 'D21' '193.29'
Set-TextValue 'E21' '  -1.76%  '
Set-TextValue 'E22' '  -2.79%  '
Set-TextValue 'D23' '6.104'
Set-TextValue 'E23' '  -4.85%  '
Set-TextValue 'E24' '  +0.89%  '
Set-TextValue 'D25' '138.41'
Set-TextValue 'E25' '  -3.35%  '
Set-TextValue 'E26' '  -2.81%  '
Set-TextValue 'D27' '7.236'
Set-TextValue 'E27' '  -3.59%  '
Set-TextValue 'D28' '16.22'
Set-TextValue 'E28' '  -0.50%  '
Set-TextValue 'D29' '1.428'
Set-TextValue 'E29' '  +0.55%  '
Set-TextValue 'D30' '0.06007'
Set-TextValue 'E30' '  -3.04%  '
Set-TextValue 'E31' '  +0.18%  '
Set-TextValue 'D32' '3.565'
Set-TextValue 'E32' '  -1.55%  '
Set-TextValue 'D33' '3.341'
Set-TextValue 'E33' '  -4.01%  '
Set-TextValue 'E34' '  -3.69%  '
Set-TextValue 'D35' '0.9826'
Set-TextValue 'E35' '  -3.80%  '
Set-TextValue 'D37' '2.781'
Set-TextValue 'E37' '  -1.00%  '
Set-TextValue 'D38' '0.5934'
Set-TextValue 'E38' '  +3.21%  '
Set-TextValue 'D39' '0.01595'
Set-TextValue 'E39' '  -3.64%  '
Set-TextValue 'D40' '5.959'
Set-TextValue 'E40' '  -0.33%  '
Set-TextValue 'D41' '0.8644'
Set-TextValue 'E41' '  -0.47%  '
Set-TextValue 'E42' '  +0.41%  '
Set-TextValue 'D43' '1.039.88'
Set-TextValue 'E43' '  -3.25%  '
Set-TextValue 'D44' '99.81'
Set-TextValue 'E44' '  -0.65%  '
Set-TextValue 'D45' '1.800.20'
Set-TextValue 'E45' '  -2.22%  '
Set-TextValue 'E46' '  +3.01%  '
Set-TextValue 'D47' '57.23'
Set-TextValue 'E47' '  -0.62%  '
Set-TextValue 'D48' '1.008'
Set-TextValue 'E48' '  +0.56%  '
Set-TextValue 'D49' '8.132'
Set-TextValue 'E49' '  -0.89%  '
Set-TextValue 'D50' '0.05180'
Set-TextValue 'E50' '  -0.74%  '
Set-TextValue 'E51' '  +3.03%  '
